$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("weibull")
$ws1.Range("B2").Value = -2.85084742220134
$ws1.Range("C2").Value = 0.373118632041817
$ws1.Range("B3").Value = 0.191139707450786
$ws1.Range("C3").Value = 0.184639301058125

$ws2 = $wb.Worksheets.Item("lognormal")
$ws2.Range("B2").Value = 2.02696016167979
$ws2.Range("C2").Value = 0.409698626967328
$ws2.Range("B3").Value = -0.968032021305164
$ws2.Range("C3").Value = 0.161971004464799

$ws3 = $wb.Worksheets.Item("llogis")
$ws3.Range("B2").Value = -2.12000443797963
$ws3.Range("C2").Value = 0.195730672493874
$ws3.Range("B3").Value = 0.506820763389294
$ws3.Range("C3").Value = 0.170291094817232

$ws4 = $wb.Worksheets.Item("gompertz")
$ws4.Range("B2").Value = -2.61749017114501
$ws4.Range("C2").Value = 0.28894923044443
$ws4.Range("B3").Value = 0.0123792653404658
$ws4.Range("C3").Value = 0.0306305115929061

$ws6 = $wb.Worksheets.Item("weibull cov")
$ws6.Range("A2").Value = 0.139217513576757
$ws6.Range("B2").Value = -0.0563951175806997
$ws6.Range("A3").Value = -0.0563951175806997
$ws6.Range("B3").Value = 0.034091671495233

$ws7 = $wb.Worksheets.Item("lognormal cov")
$ws7.Range("A2").Value = 0.167852964938914
$ws7.Range("B2").Value = -0.0587283707432745
$ws7.Range("A3").Value = -0.0587283707432745
$ws7.Range("B3").Value = 0.026234606287336

$ws8 = $wb.Worksheets.Item("llogis cov")
$ws8.Range("A2").Value = 0.0383104961549042
$ws8.Range("B2").Value = -0.00378819639480207
$ws8.Range("A3").Value = -0.00378819639480207
$ws8.Range("B3").Value = 0.0289990569740515

$ws9 = $wb.Worksheets.Item("gompertz cov")
$ws9.Range("A2").Value = 0.0834916577744283
$ws9.Range("B2").Value = -0.0063487416996439
$ws9.Range("A3").Value = -0.0063487416996439
$ws9.Range("B3").Value = 0.000938228240443155
